# Update the "kode wilayah" reference values in Sheet1, column A (rows 2-9)
# and move the active selection onto that column, per the author's
# "Sesuaikan format upload dan seeder" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$codes = @(
    "53.06.13.2001",
    "53.06.13.2002",
    "53.06.13.2003",
    "53.06.13.2004",
    "53.06.13.2005",
    "53.06.13.2006",
    "53.06.13.2007",
    "53.06.13.2008"
)

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $codes[$i]
}

# Select A2:A9 (active cell A2), matching the saved selection in the sheet.
$ws.Range("A2:A9").Select() | Out-Null
